$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 2-5) all become "b" instead of n1/n2/n10/n4
$ws.Range("A2").Value = "b"
$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "b"
$ws.Range("A5").Value = "b"

# Remove the stray E4 "paste" cell
$ws.Range("E4").ClearContents()

# Header cells fill color -> yellow (was theme color 5)
$ws.Range("A1:B1").Interior.Color = 65535

# Update selection to E15 (matches diff's sheetView selection)
$ws.Range("E15").Select()
